# Update response options for gender + license (demo-survey-short-v1)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "survey" sheet - reword the gender question label (English + Spanish)
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Range("C6").Value2 = "What best describes your gender?"
$survey.Range("D6").Value2 = "¿Cuál describe de la mejor manera su género?"

# ---------------------------------------------------------------------------
# 2. "choices" sheet
#    a) the race/ethnicity block (list hw6ou02) currently sits on rows 7-16;
#       it needs to move down one row (8-17) to make room for the extra
#       gender choice being added below.
#    b) the gender list (sj0gn93) is rewritten: "woman" moves to the top,
#       "man" becomes order 1, the old combined "non-binary / genderqueer /
#       gender non-conforming" option is replaced by two separate options,
#       "transgender" (order 2) and "non_binary" (order 3); "other" and
#       "prefer_not_to_say" shift to orders 4 and 5.
#    c) the driver's-license list (mw0ph17) loses its "prefer_not_to_say"
#       option, leaving just "yes"/"no" - this absorbs the extra row used
#       by gender so every row below stays put.
# ---------------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

# -- a) shift race/ethnicity (hw6ou02) rows 7-16 down to rows 8-17 (copy
#       bottom-up so we never clobber a row before it has been read) -------
for ($r = 16; $r -ge 7; $r--) {
    $dest = $r + 1
    $choices.Range("A$dest").Value2 = $choices.Range("A$r").Value2
    $choices.Range("B$dest").Value2 = $choices.Range("B$r").Value2
    $choices.Range("C$dest").Value2 = $choices.Range("C$r").Value2
    $choices.Range("D$dest").Value2 = $choices.Range("D$r").Value2
}

# -- b) rewrite the gender list (sj0gn93) on rows 2-7 -----------------------
$choices.Range("A2").Value2 = "sj0gn93"
$choices.Range("B2").Value2 = "woman"
$choices.Range("C2").Value2 = "Woman"
$choices.Range("D2").Value2 = "Mujer"
$choices.Range("E2").ClearContents()

$choices.Range("A3").Value2 = "sj0gn93"
$choices.Range("B3").Value2 = "man"
$choices.Range("C3").Value2 = "Man"
$choices.Range("D3").Value2 = "Hombre"
$choices.Range("E3").Value2 = 1

$choices.Range("A4").Value2 = "sj0gn93"
$choices.Range("B4").Value2 = "transgender"
$choices.Range("C4").Value2 = "Transgender"
$choices.Range("D4").Value2 = "Transgénero"
$choices.Range("E4").Value2 = 2

$choices.Range("A5").Value2 = "sj0gn93"
$choices.Range("B5").Value2 = "non_binary"
$choices.Range("C5").Value2 = "Non-binary"
$choices.Range("D5").Value2 = "No binaria"
$choices.Range("E5").Value2 = 3

$choices.Range("A6").Value2 = "sj0gn93"
$choices.Range("B6").Value2 = "other"
$choices.Range("C6").Value2 = "Other"
$choices.Range("D6").Value2 = "Otro"
$choices.Range("E6").Value2 = 4

$choices.Range("A7").Value2 = "sj0gn93"
$choices.Range("B7").Value2 = "prefer_not_to_say"
$choices.Range("C7").Value2 = "Prefer not to say"
$choices.Range("D7").Value2 = "Prefiero no decirlo"
$choices.Range("E7").Value2 = 5

# -- c) driver's-license list (mw0ph17): used to be rows 17-19
#       (yes/no/prefer_not_to_say); now it is just rows 18-19 (yes/no).
#       Rows 20 onward (the is7jb99 education list, etc.) do not move at
#       all: the one extra row gender gained above is exactly offset by
#       the one row license loses here, so nothing past row 19 changes. --
$choices.Range("A18").Value2 = "mw0ph17"
$choices.Range("B18").Value2 = "yes"
$choices.Range("C18").Value2 = "Yes"
$choices.Range("D18").Value2 = "Sí"

$choices.Range("A19").Value2 = "mw0ph17"
$choices.Range("B19").Value2 = "no"
$choices.Range("C19").Value2 = "No"
$choices.Range("D19").Value2 = "No"
$choices.Range("E19").ClearContents()
